$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their literal text representation (as in the source
# OOXML, values are stored as text, e.g. "327.97" / "6.20%") instead of being
# auto-converted to numbers/percentages by Excel when typed in.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.20%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "10.18%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.659"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "10.61%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08135"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.05%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.574"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.34%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.719"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.03%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.941"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.72%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.943"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.27%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9439"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.40%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1308"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "14.88%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.2001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "7.32%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09286"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.54%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03476"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.19%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09632"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.89%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001320"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.93%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006338"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.89%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.376"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.05%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3534"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.53%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.691"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "21.89%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1436"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.13%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2445"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.59%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04432"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.06%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.14%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004377"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.61%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001192"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.90%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003983"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "37.08%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02528"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "18.93%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05280"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "7.43%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007594"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.06%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1434"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.07%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008948"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.63%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01086"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "26.23%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006752"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.46%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.37%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002876"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-12.85%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001797"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "24.29%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.37%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.37%"
